$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells remain text (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.759.34"
$ws.Range("D3").Value = "1.746.08"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "235.11"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "0.5075"
$ws.Range("E7").Value = "  +3.25%  "
$ws.Range("D8").Value = "40.42"
$ws.Range("E8").Value = "  -2.70%  "
$ws.Range("D9").Value = "0.2665"
$ws.Range("E9").Value = "  +6.90%  "
$ws.Range("D10").Value = "0.06174"
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("D11").Value = "1.754.13"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "0.06945"
$ws.Range("E12").Value = "  +1.93%  "
$ws.Range("D13").Value = "15.34"
$ws.Range("E13").Value = "  +3.91%  "
$ws.Range("E14").Value = "  +11.10%  "
$ws.Range("D15").Value = "4.465"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "77.52"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "25.779.21"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "11.57"
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("D21").Value = "0.000006638"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("D22").Value = "1.977.07"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").Value = "4.044"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").Value = "8.238"
$ws.Range("E24").Value = "  +5.08%  "
$ws.Range("D25").Value = "5.126"
$ws.Range("E25").Value = "  +2.61%  "
$ws.Range("D26").Value = "136.44"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "1.456"
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").Value = "15.04"
$ws.Range("E28").Value = "  +2.75%  "
$ws.Range("D29").Value = "1.756"
$ws.Range("E29").Value = "  -2.50%  "
$ws.Range("D30").Value = "102.44"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").Value = "0.08170"
$ws.Range("E31").Value = "  +2.17%  "
$ws.Range("D32").Value = "3.686"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").Value = "3.383"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").Value = "0.04395"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").Value = "2.650"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("D36").Value = "0.9910"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("D37").Value = "0.5996"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "2.532"
$ws.Range("E38").Value = "  -5.38%  "
$ws.Range("D39").Value = "0.01555"
$ws.Range("E39").Value = "  +4.21%  "
$ws.Range("D40").Value = "1.935"
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "101.35"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("D43").Value = "0.3808"
$ws.Range("E43").Value = "  +3.00%  "
$ws.Range("D44").Value = "0.7448"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("D45").Value = "4.876"
$ws.Range("E45").Value = "  -4.94%  "
$ws.Range("D46").Value = "0.05498"
$ws.Range("E46").Value = "  +5.52%  "
$ws.Range("D47").Value = "0.1092"
$ws.Range("E47").Value = "  +2.68%  "
$ws.Range("D48").Value = "5.915"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("D51").Value = "1.005"
$ws.Range("E51").Value = "  +0.65%  "
